$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each updated cell as plain text (NumberFormat "@" forces text
# interpretation so numeric-looking strings like "212.60" or "29.650.86"
# keep their exact original formatting/type instead of being coerced to
# a number by Excel's automatic type detection), then ClearFormats so the
# cell style index is left untouched (matches source, which has no `s` attr).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "29.665.39"
Set-TextValue $ws.Range("E2") "  +0.85%  "
Set-TextValue $ws.Range("D3") "1.614.48"
Set-TextValue $ws.Range("E3") "  +0.64%  "
Set-TextValue $ws.Range("E4") "  -0.55%  "
Set-TextValue $ws.Range("D5") "212.60"
Set-TextValue $ws.Range("E5") "  +0.16%  "
Set-TextValue $ws.Range("D6") "0.521"
Set-TextValue $ws.Range("E6") "  +0.32%  "
Set-TextValue $ws.Range("D7") "0.992"
Set-TextValue $ws.Range("E7") "  -0.63%  "
Set-TextValue $ws.Range("D8") "29.12"
Set-TextValue $ws.Range("E8") "  +9.69%  "
Set-TextValue $ws.Range("E10") "  +1.49%  "
Set-TextValue $ws.Range("D11") "0.0910"
Set-TextValue $ws.Range("D12") "1.844.83"
Set-TextValue $ws.Range("E12") "  +0.68%  "
Set-TextValue $ws.Range("D13") "1.615.36"
Set-TextValue $ws.Range("E13") "  +0.50%  "
Set-TextValue $ws.Range("E14") "  +6.90%  "
Set-TextValue $ws.Range("D15") "3.90"
Set-TextValue $ws.Range("E15") "  +5.77%  "
Set-TextValue $ws.Range("D16") "29.678.15"
Set-TextValue $ws.Range("E16") "  +0.81%  "
Set-TextValue $ws.Range("D17") "8.98"
Set-TextValue $ws.Range("E17") "  +17.44%  "
Set-TextValue $ws.Range("D18") "64.18"
Set-TextValue $ws.Range("D19") "241.01"
Set-TextValue $ws.Range("E19") "  +0.12%  "
Set-TextValue $ws.Range("D20") "0.0₃0709"
Set-TextValue $ws.Range("E20") "  +2.91%  "
Set-TextValue $ws.Range("E21") "  -0.46%  "
Set-TextValue $ws.Range("E22") "  +2.93%  "
Set-TextValue $ws.Range("D23") "9.67"
Set-TextValue $ws.Range("E23") "  +5.96%  "
Set-TextValue $ws.Range("D24") "2.12"
Set-TextValue $ws.Range("E24") "  +1.54%  "
Set-TextValue $ws.Range("D25") "156.38"
Set-TextValue $ws.Range("E25") "  +1.30%  "
Set-TextValue $ws.Range("E27") "  +1.71%  "
Set-TextValue $ws.Range("E28") "  +3.55%  "
Set-TextValue $ws.Range("E29") "  -0.49%  "
Set-TextValue $ws.Range("E30") "  +3.36%  "
Set-TextValue $ws.Range("E31") "  +2.96%  "
Set-TextValue $ws.Range("E32") "  +2.87%  "
Set-TextValue $ws.Range("E33") "  +3.84%  "
Set-TextValue $ws.Range("D34") "1.424.52"
Set-TextValue $ws.Range("E34") "  +0.80%  "
Set-TextValue $ws.Range("E35") "  +6.81%  "
Set-TextValue $ws.Range("D36") "1.04"
Set-TextValue $ws.Range("E36") "  +1.07%  "
Set-TextValue $ws.Range("D37") "2.87"
Set-TextValue $ws.Range("E37") "  +1.65%  "
Set-TextValue $ws.Range("E38") "  -1.11%  "
Set-TextValue $ws.Range("E39") "  +3.02%  "
Set-TextValue $ws.Range("D40") "0.557"
Set-TextValue $ws.Range("E40") "  +4.16%  "
Set-TextValue $ws.Range("D41") "0.0505"
Set-TextValue $ws.Range("E41") "  +3.45%  "
Set-TextValue $ws.Range("E42") "  +4.26%  "
Set-TextValue $ws.Range("E43") "  -0.86%  "
Set-TextValue $ws.Range("B44") "Aave"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D44") "69.75"
Set-TextValue $ws.Range("E44") "  +6.32%  "
Set-TextValue $ws.Range("B45") "BitcoinSV"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D45") "53.86"
Set-TextValue $ws.Range("E45") "  +2.64%  "
Set-TextValue $ws.Range("D46") "0.992"
Set-TextValue $ws.Range("E46") "  -0.61%  "
Set-TextValue $ws.Range("E47") "  +18.62%  "
Set-TextValue $ws.Range("D48") "5.43"
Set-TextValue $ws.Range("E48") "  +3.38%  "
Set-TextValue $ws.Range("D49") "1.754.23"
Set-TextValue $ws.Range("E49") "  +0.56%  "
Set-TextValue $ws.Range("D50") "87.82"
Set-TextValue $ws.Range("E50") "  +1.35%  "
Set-TextValue $ws.Range("B51") "mCoin"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
Set-TextValue $ws.Range("D51") "2.11"
Set-TextValue $ws.Range("E51") "  -1.70%  "
